$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill in the "Responsable" (B) and "estado" (C) columns for the two new
# task rows (48 and 49) that previously only had column A populated.
$ws.Range("B48").Value = "Agustina"
$ws.Range("C48").Value = "en proceso"
$ws.Range("B49").Value = "Agustina"
$ws.Range("C49").Value = "en proceso"

# Update the active cell / selection on the sheet to C50.
$ws.Range("C50").Select()
